$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (current row 3 and below shift down by one)
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly record
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44631
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100112052
$ws.Cells.Item(3, 7).Value = "Albahaca"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 12).Value = 8000
$ws.Cells.Item(3, 13).Value = 7500
$ws.Cells.Item(3, 14).Value = "$/docena de matas"
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 1250
$ws.Cells.Item(3, 17).Value = 6
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# Copy the date formatting style from row 2's date cell to the new row's date cell
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)
